# Update countries & provincias Spain
# Applies the COVID-19 country stats refresh captured in the commit:
#  - Pakistan overtakes Mexico in total cases (row swap with updated Pakistan data)
#  - Australia, Mongolia, Butan, Gambia get refreshed counts in place
#  - Fiyi/Dominica swap display order (identical underlying counts)
#  - Islas Turcas y Caicos overtakes Santa Sede (row swap, data travels with country)
#  - Islas Virgenes Britanicas overtakes Papua Nueva Guinea (row swap, data travels with country)
#  - Footer timestamp refreshed from 04:57 to 06:14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 06:14"

# Rows 16 & 17: Pakistan now ranks above Mexico, with refreshed Pakistan totals;
# Mexico's row keeps its previous totals, just shifted one row down.
$ws.Range("A16").Value = "Pakistan"
$ws.Range("B16").Value = 176617
$ws.Range("C16").Value = 4951
$ws.Range("D16").Value = 67892
$ws.Range("E16").Value = 105224
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 119
$ws.Range("H16").Value = 3501

$ws.Range("A17").Value = "Mexico"
$ws.Range("B17").Value = 175202
$ws.Range("C17").Value = 4717
$ws.Range("D17").Value = 131686
$ws.Range("E17").Value = 22735
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 387
$ws.Range("H17").Value = 20781

# Row 73: Australia refreshed counts
$ws.Range("B73").Value = 7461
$ws.Range("C73").Value = 25
$ws.Range("D73").Value = 6896
$ws.Range("E73").Value = 463

# Row 164: Mongolia refreshed counts
$ws.Range("B164").Value = 206
$ws.Range("C164").Value = 2
$ws.Range("E164").Value = 67

# Row 185: Butan refreshed counts
$ws.Range("D185").Value = 30
$ws.Range("E185").Value = 38

# Row 191: Gambia refreshed counts
$ws.Range("B191").Value = 37
$ws.Range("C191").Value = 1
$ws.Range("G191").Value = 1
$ws.Range("H191").Value = 2

# Rows 202 & 203: Fiyi and Dominica swap display order (same underlying counts)
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# Rows 208 & 209: Islas Turcas y Caicos now ranks above Santa Sede; data travels with the country
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# Rows 213 & 214: Islas Virgenes Britanicas now ranks above Papua Nueva Guinea; data travels with the country
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
